# Updated cryptos list on Tue Nov  5 13:46:59 UTC 2024 with GitHub Actions
#
# Refreshes the Price (col D) and Volume(1h) (col E) snapshot values pulled from
# coinranking.com. Rows 32/33 additionally swap places: Bittensor drops from rank
# 30 to 31 and FirstDigitalUSD moves up to rank 30, so columns B/C (name/link) and
# D/E (price/volume) are rewritten there too.
#
# Every Value write is prefixed with a literal "'" (quote-prefix) to force text
# entry -- column D stores numeric-looking strings, some using "." as a thousands
# separator (e.g. "68.782.12"), which Excel would otherwise coerce into a real
# number/date. ".Style = 'Normal'" immediately afterwards clears the quote-prefix
# cell style that the text coercion leaves behind, so formatting stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "'68.782.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.23%  "
$ws.Range("E2").Style = "Normal"

# Row 3: Ethereum
$ws.Range("D3").Value = "'2.442.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.08%  "
$ws.Range("E3").Style = "Normal"

# Row 4: TetherUSD
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"

# Row 5: BNB
$ws.Range("D5").Value = "'560.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.07%  "
$ws.Range("E5").Style = "Normal"

# Row 6: Solana
$ws.Range("D6").Value = "'163.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.33%  "
$ws.Range("E6").Style = "Normal"

# Row 7: USDC
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("E7").Style = "Normal"

# Row 8: XRP
$ws.Range("E8").Value = "'  -0.77%  "
$ws.Range("E8").Style = "Normal"

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.171"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +8.26%  "
$ws.Range("E9").Style = "Normal"

# Row 10: TRON
$ws.Range("E10").Value = "'  -2.19%  "
$ws.Range("E10").Style = "Normal"

# Row 11: Cardano
$ws.Range("E11").Value = "'  +0.28%  "
$ws.Range("E11").Style = "Normal"

# Row 12: Toncoin
$ws.Range("E12").Value = "'  -5.04%  "
$ws.Range("E12").Style = "Normal"

# Row 13: ShibaInu
$ws.Range("E13").Value = "'  +4.71%  "
$ws.Range("E13").Style = "Normal"

# Row 14: WrappedBTC
$ws.Range("D14").Value = "'68.658.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.27%  "
$ws.Range("E14").Style = "Normal"

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'2.890.78"
$ws.Range("D15").Style = "Normal"

# Row 16: Avalanche
$ws.Range("E16").Value = "'  -0.83%  "
$ws.Range("E16").Style = "Normal"

# Row 17: WrappedEther
$ws.Range("D17").Value = "'2.443.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.75%  "
$ws.Range("E17").Style = "Normal"

# Row 18: Chainlink
$ws.Range("D18").Value = "'10.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.15%  "
$ws.Range("E18").Style = "Normal"

# Row 19: BitcoinCash
$ws.Range("D19").Value = "'339.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.11%  "
$ws.Range("E19").Style = "Normal"

# Row 20: Uniswap
$ws.Range("D20").Value = "'7.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.17%  "
$ws.Range("E20").Style = "Normal"

# Row 21: Polkadot
$ws.Range("E21").Value = "'  +1.17%  "
$ws.Range("E21").Style = "Normal"

# Row 22: SuiNetwork
$ws.Range("E22").Value = "'  +2.27%  "
$ws.Range("E22").Style = "Normal"

# Row 23: Dai
$ws.Range("E23").Value = "'  -0.02%  "
$ws.Range("E23").Style = "Normal"

# Row 24: Litecoin
$ws.Range("D24").Value = "'65.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.97%  "
$ws.Range("E24").Style = "Normal"

# Row 25: NEARProtocol
$ws.Range("D25").Value = "'3.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.79%  "
$ws.Range("E25").Style = "Normal"

# Row 26: WrappedeETH
$ws.Range("D26").Value = "'2.568.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.12%  "
$ws.Range("E26").Style = "Normal"

# Row 27: Aptos
$ws.Range("D27").Value = "'8.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.02%  "
$ws.Range("E27").Style = "Normal"

# Row 28: Binance-PegBSC-USD
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.09%  "
$ws.Range("E28").Style = "Normal"

# Row 29: PEPE
$ws.Range("E29").Value = "'  -0.07%  "
$ws.Range("E29").Style = "Normal"

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").Value = "'7.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.37%  "
$ws.Range("E30").Style = "Normal"

# Row 31: Fetch.AI
$ws.Range("E31").Value = "'  +4.77%  "
$ws.Range("E31").Style = "Normal"

# Row 32: Bittensor -> FirstDigitalUSD
$ws.Range("B32").Value = "'FirstDigitalUSD"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.07%  "
$ws.Range("E32").Style = "Normal"

# Row 33: FirstDigitalUSD -> Bittensor
$ws.Range("B33").Value = "'Bittensor"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'433.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.77%  "
$ws.Range("E33").Style = "Normal"

# Row 34: PancakeSwap
$ws.Range("E34").Value = "'  -1.96%  "
$ws.Range("E34").Style = "Normal"

# Row 35: Monero
$ws.Range("D35").Value = "'159.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.09%  "
$ws.Range("E35").Style = "Normal"

# Row 36: WhiteBITCoin
$ws.Range("E36").Value = "'  -0.09%  "
$ws.Range("E36").Style = "Normal"

# Row 37: USDe
$ws.Range("E37").Value = "'  +0.02%  "
$ws.Range("E37").Style = "Normal"

# Row 38: EthereumClassic
$ws.Range("E38").Value = "'  +0.59%  "
$ws.Range("E38").Style = "Normal"

# Row 39: Kaspa
$ws.Range("E39").Value = "'  -0.42%  "
$ws.Range("E39").Style = "Normal"

# Row 40: PolygonEcosystemToken
$ws.Range("E40").Value = "'  +0.86%  "
$ws.Range("E40").Style = "Normal"

# Row 41: Stacks
$ws.Range("E41").Value = "'  +2.19%  "
$ws.Range("E41").Style = "Normal"

# Row 42: RenderToken
$ws.Range("E42").Value = "'  -1.21%  "
$ws.Range("E42").Style = "Normal"

# Row 43: ImmutableX
$ws.Range("E43").Value = "'  +0.47%  "
$ws.Range("E43").Style = "Normal"

# Row 44: dogwifhat
$ws.Range("E44").Value = "'  +0.96%  "
$ws.Range("E44").Style = "Normal"

# Row 45: Filecoin
$ws.Range("E45").Value = "'  -1.44%  "
$ws.Range("E45").Style = "Normal"

# Row 46: Aave
$ws.Range("D46").Value = "'130.01"
$ws.Range("D46").Style = "Normal"

# Row 47: Cronos
$ws.Range("E47").Value = "'  -0.31%  "
$ws.Range("E47").Style = "Normal"

# Row 48: ARBITRUM
$ws.Range("E48").Value = "'  -0.20%  "
$ws.Range("E48").Style = "Normal"

# Row 49: Mantle
$ws.Range("E49").Value = "'  -1.39%  "
$ws.Range("E49").Style = "Normal"

# Row 50: Stellar
$ws.Range("D50").Value = "'0.0925"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.52%  "
$ws.Range("E50").Style = "Normal"

# Row 51: BitgetToken
$ws.Range("E51").Value = "'  +2.55%  "
$ws.Range("E51").Style = "Normal"
